$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain text (not numbers) in the source data; force
# text format while assigning so Excel does not auto-convert strings like
# "0.654" into numeric values, then restore the original cell style.

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.838.70'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  -0.03%  '
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.269.90'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  +0.36%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.654'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +1.76%  '
$ws.Range("E6").Value = '  -0.62%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.87'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("E8").Value = '  -0.01%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.449'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +4.26%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0984'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -3.92%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.91'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +0.77%  '
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.75'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +3.05%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +1.79%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.602.88'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -0.18%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.65'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -0.66%  '
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.15'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  +3.36%  '
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.844'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +2.03%  '
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.248.01'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -0.76%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.806.80'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  +0.45%  '
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0985'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  -0.32%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.95'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -0.30%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.17'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +0.96%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.23'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +0.30%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  +32.64%  '
$ws.Range("E26").Value = '  -2.27%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.30'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -0.97%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.97'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  +0.77%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '174.49'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +0.79%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.98'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +5.14%  '
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("E33").Value = '  +2.11%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.99'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +4.85%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0685'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -0.41%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.96'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  -2.32%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.71'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -4.27%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.45'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -5.91%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.31'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("E41").Value = '  +0.11%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.65'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +2.81%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.86'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.23'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -1.77%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.45'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("E48").Value = '  +6.07%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.459.03'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -1.55%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.04'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -4.36%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.32'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  -0.59%  '
